$d = $word.ActiveDocument

# The primary header contains two floating (anchored) pictures whose
# display names need to be swapped:
#   "image2.png" -> "image1.png"
#   "image1.jpg" -> "image2.jpg"
# This only renames the picture (Shape.Name / docPr "name" attribute);
# the underlying embedded image data and relationship targets (rId1 ->
# media/image2.png, rId2 -> media/image1.jpg) are left untouched.

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($hfIdx = 1; $hfIdx -le 3; $hfIdx++) {
        $hdr = $sec.Headers.Item($hfIdx)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Shapes.Count; $i++) {
                $shp = $hdr.Shapes.Item($i)
                if ($shp.Name -eq "image2.png") {
                    $shp.Name = "image1.png"
                } elseif ($shp.Name -eq "image1.jpg") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
}
